# Update the "Info" sheet so that the @base and @prefix URI/URL values are
# wrapped in angle brackets, matching the commit:
# "Updates to require angle brackets for @base and @prefix values."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Range("D1").Value = "<http://sales.data/purchases/2015>"
$ws.Range("D2").Value = "<http://sales.data/purchases#>"
$ws.Range("D3").Value = "<http://sales.data/schema#>"

# Update the sheet's current selection to reflect the cells just edited.
$ws.Range("D1:D3").Select()
